$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1) - column F ("想去人数") updates
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 351
$ws1.Cells.Item(4, 6).Value = 1280
$ws1.Cells.Item(10, 6).Value = 3527
$ws1.Cells.Item(11, 6).Value = 138
$ws1.Cells.Item(12, 6).Value = 92
$ws1.Cells.Item(13, 6).Value = 71
$ws1.Cells.Item(14, 6).Value = 47
$ws1.Cells.Item(15, 6).Value = 59
$ws1.Cells.Item(16, 6).Value = 611
$ws1.Cells.Item(17, 6).Value = 100
$ws1.Cells.Item(18, 6).Value = 764
$ws1.Cells.Item(19, 6).Value = 212
$ws1.Cells.Item(20, 6).Value = 128
$ws1.Cells.Item(22, 6).Value = 64
$ws1.Cells.Item(23, 6).Value = 69
$ws1.Cells.Item(24, 6).Value = 2711
$ws1.Cells.Item(25, 6).Value = 5212
$ws1.Cells.Item(28, 6).Value = 479
$ws1.Cells.Item(29, 6).Value = 3085
$ws1.Cells.Item(30, 6).Value = 292
$ws1.Cells.Item(31, 6).Value = 2261
$ws1.Cells.Item(35, 6).Value = 128
$ws1.Cells.Item(36, 6).Value = 181
$ws1.Cells.Item(38, 6).Value = 33
$ws1.Cells.Item(40, 6).Value = 810
$ws1.Cells.Item(42, 6).Value = 5
$ws1.Cells.Item(44, 6).Value = 40
$ws1.Cells.Item(45, 6).Value = 489

# Sheet "演出" (index 2) - column F ("想去人数") updates
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(2, 6).Value = 75

# Sheet "全部类型" (index 4) - column F ("想去人数") updates
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 351
$ws4.Cells.Item(4, 6).Value = 1280
$ws4.Cells.Item(10, 6).Value = 3527
$ws4.Cells.Item(11, 6).Value = 138
$ws4.Cells.Item(12, 6).Value = 92
$ws4.Cells.Item(13, 6).Value = 71
$ws4.Cells.Item(14, 6).Value = 75
$ws4.Cells.Item(15, 6).Value = 47
$ws4.Cells.Item(16, 6).Value = 59
$ws4.Cells.Item(17, 6).Value = 611
$ws4.Cells.Item(18, 6).Value = 100
$ws4.Cells.Item(19, 6).Value = 764
$ws4.Cells.Item(20, 6).Value = 212
$ws4.Cells.Item(21, 6).Value = 128
$ws4.Cells.Item(23, 6).Value = 64
$ws4.Cells.Item(24, 6).Value = 69
$ws4.Cells.Item(25, 6).Value = 2711
$ws4.Cells.Item(26, 6).Value = 5212
$ws4.Cells.Item(29, 6).Value = 479
$ws4.Cells.Item(30, 6).Value = 3085
$ws4.Cells.Item(31, 6).Value = 292
$ws4.Cells.Item(32, 6).Value = 2261
$ws4.Cells.Item(36, 6).Value = 128
$ws4.Cells.Item(37, 6).Value = 181
$ws4.Cells.Item(39, 6).Value = 33
$ws4.Cells.Item(41, 6).Value = 810
$ws4.Cells.Item(43, 6).Value = 5
$ws4.Cells.Item(45, 6).Value = 40
$ws4.Cells.Item(46, 6).Value = 489

Write-Output "Done updating F column values."
